# Updated cryptos list with latest price/volume data.
# Applies the per-cell text updates while forcing the cells to remain
# plain text (matching their original inlineStr/text storage) and
# keeping them on the workbook's default (unstyled) cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "532.51")
    # are not reinterpreted as numbers by Excel's input parser.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    # Cells in this sheet use the default/general style (no explicit
    # number format) -- restore that so we don't leave a stray
    # text-format style applied to the cell.
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '57.797.48'
Set-TextValue 'E2' '  +1.42%  '
Set-TextValue 'D3' '3.116.02'
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '532.51'
Set-TextValue 'D6' '138.24'
Set-TextValue 'E6' '  +1.56%  '
Set-TextValue 'E7' '  -0.16%  '
Set-TextValue 'D8' '0.499'
Set-TextValue 'E8' '  +10.76%  '
Set-TextValue 'D9' '7.35'
Set-TextValue 'E9' '  +0.30%  '
Set-TextValue 'E10' '  +1.42%  '
Set-TextValue 'D11' '0.413'
Set-TextValue 'E11' '  +3.94%  '
Set-TextValue 'E12' '  +3.50%  '
Set-TextValue 'D13' '3.652.29'
Set-TextValue 'E13' '  +1.13%  '
Set-TextValue 'D14' '25.66'
Set-TextValue 'E14' '  +1.50%  '
Set-TextValue 'E15' '  +3.61%  '
Set-TextValue 'D16' '57.901.90'
Set-TextValue 'E16' '  +1.41%  '
Set-TextValue 'D17' '3.115.45'
Set-TextValue 'E17' '  +1.19%  '
Set-TextValue 'D18' '6.14'
Set-TextValue 'E18' '  +4.41%  '
Set-TextValue 'E19' '  +2.65%  '
Set-TextValue 'D20' '8.11'
Set-TextValue 'E20' '  +3.25%  '
Set-TextValue 'D21' '375.10'
Set-TextValue 'E21' '  +8.05%  '
Set-TextValue 'E22' '  +0.07%  '
Set-TextValue 'D23' '5.74'
Set-TextValue 'E23' '  -1.81%  '
Set-TextValue 'D24' '69.51'
Set-TextValue 'E24' '  +1.77%  '
Set-TextValue 'E25' '  +1.93%  '
Set-TextValue 'E26' '  -0.13%  '
Set-TextValue 'E27' '  -0.14%  '
Set-TextValue 'E28' '  +1.89%  '
Set-TextValue 'E29' '  +4.96%  '
Set-TextValue 'E30' '  +4.39%  '
Set-TextValue 'E31' '  +0.09%  '
Set-TextValue 'D32' '21.52'
Set-TextValue 'E32' '  +3.51%  '
Set-TextValue 'D33' '5.14'
Set-TextValue 'E33' '  +4.44%  '
Set-TextValue 'E34' '  +2.93%  '
Set-TextValue 'D35' '160.54'
Set-TextValue 'E35' '  +0.86%  '
Set-TextValue 'E36' '  +2.97%  '
Set-TextValue 'E37' '  +6.49%  '
Set-TextValue 'D38' '25.52'
Set-TextValue 'E38' '  -1.31%  '
Set-TextValue 'E39' '  +3.86%  '
Set-TextValue 'D40' '0.0670'
Set-TextValue 'E40' '  +2.79%  '
Set-TextValue 'D41' '2.561.48'
Set-TextValue 'E41' '  +7.15%  '
Set-TextValue 'E42' '  +3.62%  '
Set-TextValue 'D43' '38.61'
Set-TextValue 'E43' '  +5.47%  '
Set-TextValue 'D44' '0.697'
Set-TextValue 'E44' '  +0.93%  '
Set-TextValue 'E45' '  +2.46%  '
Set-TextValue 'E46' '  -0.01%  '
Set-TextValue 'E47' '  +2.17%  '
Set-TextValue 'E48' '  +3.49%  '
Set-TextValue 'D49' '19.93'
Set-TextValue 'E49' '  +1.61%  '
Set-TextValue 'D50' '0.0952'
Set-TextValue 'E50' '  +6.95%  '
Set-TextValue 'E51' '  -0.89%  '
